# "llama models FT corrected"
#
# The "Llama-3-8b" row (row 4) on every sheet was missing its "X" mark in
# column B ("AirlinesCustomerSatisfaction" / first data column) even
# though the row's Total already counted it - add the missing mark back
# on all three sheets, then leave each sheet's selection where the
# author left it while reviewing the fix.

$wb = $excel.ActiveWorkbook

# --- Sheet1 --------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("B4").Value = "X"

# --- Sheet2 (stays the active/visible tab) --------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("B4").Value = "X"

# --- Sheet3 ----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("B4").Value = "X"

# Restore each sheet's reviewing cursor position.
$ws1.Range("C6").Select() | Out-Null
$ws3.Range("C8").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("F6").Select() | Out-Null
